# buffer grilla zona incendios
# Applies the edits observed in the target diff:
#  1. Add a new comment/observation in J3 (row for the "Satellite-Based Fire
#     Progression Mapping" paper).
#  2. Clear the "Baja importancia" note in I4 (Spatial variability paper row).
#  3. Remove the "Le va pésimo" note that lived in I5 (Next Day Wildfire
#     Spread paper row).
#  4. Add a new paper entry in row 13: "Comparing Next-Day Wildfire
#     Predictability of MODIS and VIIRS Satellite Data", with its arxiv link
#     hyperlinked in B13 (matching the style used by the other paper rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1. New note in J3
$ws.Range("J3").Value = "Según los puntos donde sabemos que se quemó algo en un día dado, estimamos por interpolación cuándo (qué día) se quemaron las demás celdas`""

# 2. Clear I4 ("Baja importancia") but keep its existing style/format
$ws.Range("I4").ClearContents()

# 3. Clear I5 ("Le va pésimo") entirely
$ws.Range("I5").ClearContents()

# 4. New paper row 13 (link entered before title, matching shared-string order)
$ws.Range("B13").Value = "https://arxiv.org/pdf/2503.08580"
$ws.Hyperlinks.Add($ws.Range("B13"), "https://arxiv.org/pdf/2503.08580") | Out-Null
$ws.Range("B13").Style = $ws.Range("B9").Style

$ws.Range("A13").Value = "Comparing Next-Day Wildfire Predictability of`nMODIS and VIIRS Satellite Data"
$ws.Range("A13").WrapText = $true
